$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AA1")
$ws.Select()

# The two calibration blocks (curvature x/y pairs) had their x and y
# columns (A and B) swapped.
$rangeTop = $ws.Range("A4:B8")
$rangeBottom = $ws.Range("A10:B14")

$valsTop = $rangeTop.Value2
$valsBottom = $rangeBottom.Value2

$rangeTop.Value2 = $valsBottom
$rangeBottom.Value2 = $valsTop

# Leave the selection where the author ended up after editing.
$ws.Range("A3:B14").Select()
